$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-09"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 10-09)"

# Update the November (row 11) and Total (row 14) figures for the "2022 (through ...)" column (I)
$ws.Range("I11").Value = 30
$ws.Range("I14").Value = 1308
